# Week 15 logged / Week 16 simulated update for Chargers Players Data.xlsx
#
# Sheet1 "Rushing": a new row (G.Nabers) is inserted at row 7, pushing
# M.Williams / J.Palmer / J.Guyton / A.Roberts down by one row, and all
# rows' rushing totals (1DATT/2DATT/3DATT/RZATT) get updated.
#
# Sheet2 "Receiving": no rows inserted, only receiving totals updated for
# most players (A.Ekeler, J.Jackson, J.Kelley, K.Allen, M.Williams,
# J.Palmer, J.Guyton, K.Hill(name only), J.Cook, D.Parham, S.Anderson,
# T.McKitty).
#
# Active sheet/selection also flips from Receiving -> Rushing.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Rushing")
$ws2 = $wb.Worksheets.Item("Receiving")

# ---------------------------------------------------------------------
# Sheet1 (Rushing): insert new row for G.Nabers at row 7
# ---------------------------------------------------------------------
$ws1.Rows("7:7").Insert() | Out-Null

# Copy formatting (bold/border/centered style) from the row above so the
# new "week index" cell in column A matches the rest of the table.
$ws1.Range("A6").Copy() | Out-Null
$ws1.Range("A7").PasteSpecial(-4122) | Out-Null

# New row 7: G.Nabers
$ws1.Range("A7").Value = 5
$ws1.Range("B7").Value = "G.Nabers"
$ws1.Range("C7").Value = 1
$ws1.Range("D7").Value = 0
$ws1.Range("E7").Value = 0
$ws1.Range("F7").Value = 1

# Row 2: J.Herbert
$ws1.Range("C2").Value = 14
$ws1.Range("D2").Value = 8
$ws1.Range("E2").Value = 19
$ws1.Range("F2").Value = 12

# Row 3: A.Ekeler
$ws1.Range("C3").Value = 0
$ws1.Range("D3").Value = 0
$ws1.Range("E3").Value = 0
$ws1.Range("F3").Value = 0

# Row 4: J.Jackson
$ws1.Range("C4").Value = 75
$ws1.Range("D4").Value = 35
$ws1.Range("E4").Value = 6
$ws1.Range("F4").Value = 25

# Row 5: J.Kelley
$ws1.Range("C5").Value = 30
$ws1.Range("D5").Value = 22
$ws1.Range("E5").Value = 4
$ws1.Range("F5").Value = 15

# Row 6: L.Rountree - unchanged

# Row 8 (was row 7 pre-insert): M.Williams - update week index + stats
$ws1.Range("A8").Value = 6
$ws1.Range("C8").Value = 0
$ws1.Range("D8").Value = 0
$ws1.Range("E8").Value = 1
$ws1.Range("F8").Value = 0

# Row 9 (was row 8 pre-insert): J.Palmer - update week index + stats
$ws1.Range("A9").Value = 7
$ws1.Range("C9").Value = 1
$ws1.Range("D9").Value = 0
$ws1.Range("E9").Value = 0
$ws1.Range("F9").Value = 0

# Row 10 (was row 9 pre-insert): J.Guyton - update week index + stats
$ws1.Range("A10").Value = 8
$ws1.Range("C10").Value = 4
$ws1.Range("D10").Value = 2
$ws1.Range("E10").Value = 0
$ws1.Range("F10").Value = 0

# Row 11 (was row 10 pre-insert): A.Roberts - update week index + stats
$ws1.Range("A11").Value = 9
$ws1.Range("C11").Value = 3
$ws1.Range("D11").Value = 0
$ws1.Range("E11").Value = 0
$ws1.Range("F11").Value = 1

# ---------------------------------------------------------------------
# Sheet2 (Receiving): update receiving totals (no rows inserted)
# ---------------------------------------------------------------------

# Row 2: A.Ekeler
$ws2.Range("C2").Value = 0
$ws2.Range("D2").Value = 0
$ws2.Range("E2").Value = 0
$ws2.Range("F2").Value = 0
$ws2.Range("G2").Value = 0
$ws2.Range("H2").Value = 0

# Row 3: J.Jackson
$ws2.Range("C3").Value = 30
$ws2.Range("D3").Value = 20
$ws2.Range("E3").Value = 2
$ws2.Range("F3").Value = 1

# Row 4: J.Kelley
$ws2.Range("C4").Value = 15
$ws2.Range("D4").Value = 14
$ws2.Range("E4").Value = 2
$ws2.Range("F4").Value = 1

# Row 5: L.Rountree - unchanged

# Row 6: K.Allen
$ws2.Range("C6").Value = 108
$ws2.Range("D6").Value = 75
$ws2.Range("E6").Value = 19
$ws2.Range("F6").Value = 12
$ws2.Range("G6").Value = 20
$ws2.Range("H6").Value = 13

# Row 7: M.Williams
$ws2.Range("C7").Value = 68
$ws2.Range("D7").Value = 41
$ws2.Range("E7").Value = 26
$ws2.Range("F7").Value = 12
$ws2.Range("G7").Value = 18

# Row 8: J.Palmer
$ws2.Range("C8").Value = 20
$ws2.Range("D8").Value = 15

# Row 9: J.Guyton
$ws2.Range("C9").Value = 24
$ws2.Range("D9").Value = 14
$ws2.Range("G9").Value = 5
$ws2.Range("H9").Value = 2

# Row 10: K.Hill - stats unchanged

# Row 11: J.Cook
$ws2.Range("C11").Value = 51
$ws2.Range("D11").Value = 31
$ws2.Range("E11").Value = 11
$ws2.Range("G11").Value = 7

# Row 12: D.Parham
$ws2.Range("C12").Value = 20
$ws2.Range("G12").Value = 5

# Row 13: S.Anderson
$ws2.Range("C13").Value = 12
$ws2.Range("D13").Value = 10

# Row 14: T.McKitty
$ws2.Range("C14").Value = 3
$ws2.Range("D14").Value = 2

# ---------------------------------------------------------------------
# Active sheet / selection: Rushing becomes the active tab (was Receiving)
# ---------------------------------------------------------------------
$ws2.Range("G5").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("F4").Select() | Out-Null
